$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "DES_ORGAO" abbreviation cells in column A, added for rows 10, 13, 14, 15.
# Entered in this order so the new shared-string table entries land in the
# same sequence as the target workbook (SAP/MAPA, SPRC/MAPA, SMC/MAPA, SIMS/MAPA).
$ws.Range("A14").Value = "SAP/MAPA"
$ws.Range("A15").Value = "SPRC/MAPA"
$ws.Range("A10").Value = "SMC/MAPA"
$ws.Range("A13").Value = "SIMS/MAPA"

# Match the author's final cell selection.
$ws.Range("B13").Select()
